$d = $word.ActiveDocument

# The "-robber" paragraph currently carries the _GoBack bookmark at the end
# of its text. We need to: split it off into its own paragraph (no bookmark),
# add a blank paragraph, then "LOC:", "Trent: 1,744" (now carrying the
# bookmark), "Steven: 1,222" and "CJ: 1,791".

# Remove the existing _GoBack bookmark; we'll recreate it on the new
# "Trent: 1,744" paragraph below.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Locate the "-robber" paragraph by content rather than a hard-coded index.
$robberRange = $d.Content
$robberRange.Find.Execute("-robber") | Out-Null
$prefixRange = $d.Range(0, $robberRange.End)
$robberIndex = $prefixRange.Paragraphs.Count
$robberPara = $d.Paragraphs.Item($robberIndex)

# A one-off marker character so we can later find the exact insertion point
# (end of "Trent: 1,744") with an unmutated Range - Bookmarks.Add only keeps
# precise bounds when given a Range fresh off Find.Execute.
$marker = [char]1

$robberPara.Range.InsertAfter("`r`rLOC:`rTrent: 1,744" + $marker + "`rSteven: 1,222`rCJ: 1,791")

# Re-find the marker and bookmark that exact spot.
$markerRange = $d.Content
$markerRange.Find.Execute($marker) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRange)

# Remove the marker character itself, collapsing the new bookmark to a
# zero-length point right after "Trent: 1,744" (matching the original
# placement relative to "-robber").
$markerRange.Text = ""
